$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("M2").Value = 27.3

# --- Row 3 ---
$ws.Range("F3").Value = 1.3
$ws.Range("M3").Value = 8.3699999999999992

# --- Row 4 ---
$ws.Range("F4").Value = 4.3
$ws.Range("M4").Value = 14.99

# --- Row 5 ---
$ws.Range("F5").Value = 32
$ws.Range("M5").Value = 60.69

# --- Row 6 ---
# Move the existing E6 formula into F6, and put a fresh value into E6.
$ws.Range("F6").Formula = "=51178+14426+3116+440+35+1+1"
$ws.Range("E6").Value = 72000
$ws.Range("M6").Formula = "=8591+1106+203+33+6+65"

# --- Row 7 ---
$ws.Range("F7").Value = 1.3
$ws.Range("M7").Value = 135.81

# --- Row 8 ---
$ws.Range("F8").Value = 1.3
$ws.Range("M8").Value = 9.3699999999999992

# --- Row 9 ---
$ws.Range("F9").Formula = "=13436+1133+58+2+2"

# --- Row 10 ---
$ws.Range("F10").Formula = "=5726+556+36+2+2"

# --- Row 11 ---
$ws.Range("F11").Formula = "=13852+3644+800+129+14+1+2"

# --- Row 12 ---
$ws.Range("F12").Value = 35973
$ws.Range("L12").Value = 2316

# --- Row 13 ---
# Move the existing E13 formula into F13, and put a fresh value into E13.
$ws.Range("F13").Formula = "=7132+269+5+2"
$ws.Range("E13").Value = 36000

# --- Row 14 ---
# Move the existing E14 value (formula result, entered as plain number) into F14.
$ws.Range("F14").Value = 2680
$ws.Range("E14").Value = 36000

# --- Row 15 ---
# Replace the old E15 formula with a plain value, add a new formula in F15.
$ws.Range("E15").Value = 36000
$ws.Range("F15").Formula = "=30590+4706+534+42+2"

# --- Row 16 ---
# Replace the old E16 formula with a plain value, add a new formula in F16.
$ws.Range("E16").Value = 36000
$ws.Range("F16").Formula = "=10943+7877+5251+3042+1511+625+209+51+8+1+1"

# --- Row 24 ---
$ws.Range("F24").Value = 4470

# --- Selection / view state ---
$ws.Range("M6").Select()
